$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-6 from
# 2023-09-15 (serial 45184) to 2023-09-16 (serial 45185).
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45185
}
